$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as literal text
# (matches the source data which stores these as plain strings, not numbers).

$ws.Range("D2").Value = '30.729.82'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '2.124.91'
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.016'
$ws.Range("E4").Value = '  +1.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.15'
$ws.Range("E5").Value = '  +1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.014'
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4567'
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.18'
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09121'
$ws.Range("E10").Value = '  +2.00%  '
$ws.Range("E11").Value = '  +1.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.51'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '2.135.26'
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.867'
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.124'
$ws.Range("E15").Value = '  +4.76%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '97.64'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001173'
$ws.Range("E17").Value = '  +4.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.015'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06716'
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.61'
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.498'
$ws.Range("E22").Value = '  +3.19%  '
$ws.Range("D23").Value = '30.794.84'
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.02'
$ws.Range("E24").Value = '  +5.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.384'
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("D26").Value = '2.360.35'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.49'
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '165.87'
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.559'
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.06'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.211'
$ws.Range("E31").Value = '  +0.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1080'
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.651'
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.417'
$ws.Range("E34").Value = '  +4.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.960'
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.61'
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.944'
$ws.Range("E37").Value = '  +8.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02672'
$ws.Range("E38").Value = '  +4.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06899'
$ws.Range("E39").Value = '  +1.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2330'
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.65'
$ws.Range("E41").Value = '  -1.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6941'
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.268'
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.04'
$ws.Range("E44").Value = '  +7.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6496'
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.314'
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000369'
$ws.Range("E47").Value = '  +15.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.707'
$ws.Range("E48").Value = '  +1.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.261'
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.67'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07329'
$ws.Range("E51").Value = '  +4.01%  '
